$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values taken from the refreshed coinranking.com snapshot.
# Column D price strings that look like plain numbers get a leading
# apostrophe so Excel keeps them as text (matching the source inlineStr
# cells) instead of auto-converting them to numeric values.

$ws.Range("D2").Value = "'37.159.41"
$ws.Range("D3").Value = "'2.050.19"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'248.46"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'57.08"
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "'16.20"
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("D13").Value = "'0.906"
$ws.Range("E13").Value = "  +11.93%  "
$ws.Range("D14").Value = "'2.346.90"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").Value = "'5.75"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").Value = "'2.051.09"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").Value = "'18.89"
$ws.Range("E17").Value = "  +14.01%  "
$ws.Range("D18").Value = "'37.167.48"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").Value = "'74.61"
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("D20").Value = "'0.0₃0898"
$ws.Range("E20").Value = "  -2.80%  "
$ws.Range("D21").Value = "'5.46"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "'237.17"
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +3.07%  "
$ws.Range("D25").Value = "'9.65"
$ws.Range("E25").Value = "  +3.16%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'170.17"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "'2.17"
$ws.Range("E27").Value = "  -4.80%  "
$ws.Range("D28").Value = "'20.20"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("E29").Value = "  -0.93%  "
$ws.Range("D30").Value = "'5.03"
$ws.Range("E30").Value = "  +5.13%  "
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").Value = "'0.0623"
$ws.Range("D33").Value = "'4.53"
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("D34").Value = "'0.0884"
$ws.Range("E34").Value = "  -3.10%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'2.27"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("E37").Value = "  +1.08%  "
$ws.Range("D38").Value = "'1.35"
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("D39").Value = "'5.34"
$ws.Range("E39").Value = "  +14.50%  "
$ws.Range("E40").Value = "  +8.26%  "
$ws.Range("D41").Value = "'0.0983"
$ws.Range("E41").Value = "  -14.48%  "
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("D45").Value = "'96.19"
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("E46").Value = "  -1.71%  "
$ws.Range("D47").Value = "'1.274.05"
$ws.Range("E47").Value = "  -1.65%  "
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("D49").Value = "'6.83"
$ws.Range("E49").Value = "  -1.39%  "
$ws.Range("D50").Value = "'2.233.91"
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("D51").Value = "'44.48"
$ws.Range("E51").Value = "  +0.56%  "
